$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 233 (weekly price record), pushing the
# existing rows 233:248 down to 234:249.
$ws.Rows.Item(233).Insert()

# Fill the newly inserted row with the new weekly record.
$ws.Cells.Item(233, 1).Value = 5
$ws.Cells.Item(233, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(233, 3).Value = 'Maule'
$ws.Cells.Item(233, 4).Value = 44610
$ws.Cells.Item(233, 5).Value = 7
$ws.Cells.Item(233, 6).Value = 100112003
$ws.Cells.Item(233, 7).Value = 'Ajo'
$ws.Cells.Item(233, 8).Value = 'Chino'
$ws.Cells.Item(233, 9).Value = 'Primera'
$ws.Cells.Item(233, 10).Value = 200
$ws.Cells.Item(233, 11).Value = 20000
$ws.Cells.Item(233, 12).Value = 20000
$ws.Cells.Item(233, 13).Value = 20000
$ws.Cells.Item(233, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(233, 15).Value = 'China'
$ws.Cells.Item(233, 16).Value = 2000
$ws.Cells.Item(233, 17).Value = 10
$ws.Cells.Item(233, 18).Value = 'Hortaliza'
